$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before the existing row 475, shifting the rest of the
# table (old rows 475-502) down to 478-505.
$ws.Range("A475:T477").EntireRow.Insert()

# Common (unchanged across the three new rows) column values, copied from
# the surrounding records in this weekly price table.
$mercadoId = 11
$mercado   = "Vega Monumental Concepción"
$region    = "Bíobío"
$codreg    = 8
$tipo      = "Fruta"
$productoId = 100108
$producto  = "Tropicales y subtropicales"
$categoriaId = 100108006
$categoria = "Plátano"
$variedad  = "Sin especificar"
$unidad    = "$/caja 20 kilos"
$origen    = "Ecuador"
$kgUnidad  = 20

# New rows 475-477: date 44714 ("02/06/2022"), one row per calidad.
$rows = @(
    @{ Row = 475; Calidad = "Maduro";           Volumen = 200; Min = 12000; Max = 12000; Prom = 12000; PrecioKg = 600 },
    @{ Row = 476; Calidad = "Pintón";            Volumen = 400; Min = 13000; Max = 13000; Prom = 13000; PrecioKg = 650 },
    @{ Row = 477; Calidad = "Primera Pintón";    Volumen = 400; Min = 15000; Max = 15000; Prom = 15000; PrecioKg = 750 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value  = $mercadoId
    $ws.Cells.Item($row, 2).Value  = $mercado
    $ws.Cells.Item($row, 3).Value  = $region
    $ws.Cells.Item($row, 4).Value  = 44714
    $ws.Cells.Item($row, 5).Value  = $codreg
    $ws.Cells.Item($row, 6).Value  = $tipo
    $ws.Cells.Item($row, 7).Value  = $productoId
    $ws.Cells.Item($row, 8).Value  = $producto
    $ws.Cells.Item($row, 9).Value  = $categoriaId
    $ws.Cells.Item($row, 10).Value = $categoria
    $ws.Cells.Item($row, 11).Value = $variedad
    $ws.Cells.Item($row, 12).Value = $r.Calidad
    $ws.Cells.Item($row, 13).Value = $r.Volumen
    $ws.Cells.Item($row, 14).Value = $r.Min
    $ws.Cells.Item($row, 15).Value = $r.Max
    $ws.Cells.Item($row, 16).Value = $r.Prom
    $ws.Cells.Item($row, 17).Value = $unidad
    $ws.Cells.Item($row, 18).Value = $origen
    $ws.Cells.Item($row, 19).Value = $r.PrecioKg
    $ws.Cells.Item($row, 20).Value = $kgUnidad
}
